# TMTT0046984_VerificationOfEngagementCompanyRoundtripFlagFunctionalityOnTheBuysideDeals
# - Mid - 5th May 2025
#
# Changes applied:
#  1. Users sheet, cell B2: "Gemma Hardy" -> "Brian Miller".
#  2. The "AddContact" sheet is no longer the active/selected tab.
#  3. The "Users" sheet becomes the active/selected tab, with C8:C9
#     selected (previously C12 was selected and it was not the active tab).

$wb = $excel.ActiveWorkbook

$usersSheet = $wb.Worksheets.Item("Users")
$addContactSheet = $wb.Worksheets.Item("AddContact")

# 1. Fix the name in the Users table.
$usersSheet.Range("B2").Value = "Brian Miller"

# 2. Move off of AddContact ...
$addContactSheet.Select() | Out-Null

# 3. ... and onto Users, leaving C8:C9 as the selected range there.
$usersSheet.Select() | Out-Null
$usersSheet.Range("C8:C9").Select() | Out-Null
